$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: LinearRegression - update single score precision
$ws.Range("B2").Value = 0.8091271741280979

# Row 3: RandomForestRegressor - updated scores
$ws.Range("B3").Value = 0.8660749081655226
$ws.Range("C3").Value = 0.8974191652456925
$ws.Range("D3").Value = 0.8640674321421672

# Row 4: model changed from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.8338498049065466
$ws.Range("C4").Value = 0.8316492258835613
$ws.Range("D4").Value = 0.8434017444386065

# Row 5: model changed from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8698051527777229
$ws.Range("C5").Value = 0.9053963869124866
$ws.Range("D5").Value = 0.8606551045605472
